# Add a new "how to demo" column (E) with per-row demo notes, and fix the
# typo in the "初始时间估算" header (was "初试时间估算").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: 初试时间估算 -> 初始时间估算
$ws.Range("C1").Value = "初始时间估算"

# New column header
$ws.Range("E1").Value = "how to demo"

# Per-row "how to demo" notes
$ws.Range("E2").Value = "用户打开APP，主界面有点击拍照按钮，拍照上传后显示食物的种类，热量"
$ws.Range("E3").Value = "生成每周的食物热量报告并推送给用户，用户可以选择查看或不查看"
$ws.Range("E4").Value = "用户将自己当前的状态改为减肥状态，点击查看建议，就可以看到系统提供的又营养又瘦身的建议"
$ws.Range("E5").Value = "当系统检查到用户摄入营养不均衡时，会提醒用户健康饮食 ，用户若不需要可以关闭"
$ws.Range("E6").Value = "用户点击“帮我想想”模块，就会为用户根据自己的饮食爱好提供健康的饮食方案"
$ws.Range("E7").Value = "用户可以在设置中找到更多建议选项，为产品提供意见"
$ws.Range("E8").Value = "用户点开“更多知识”选项，会出现关于食物的文章"
$ws.Range("E9").Value = "用户在注册账号的时候可以选择通过QQ登录，会跳转到QQ界面"
$ws.Range("E10").Value = "主界面有一个搜索框，输入食物显示热量"
$ws.Range("E11").Value = "用户可以选择“推荐”模块，会根据用户喜好推荐周围的食物"
$ws.Range("E12").Value = "在“我的”模块中，有联系人选项，可以同步查询通讯录中的好友"

# Column E width, matching the new layout
$ws.Columns.Item(5).ColumnWidth = 59

# Row 6 grows to fit the wrapped demo text
$ws.Rows.Item(6).RowHeight = 27.6

# Selection ends on the last edited cell
$ws.Range("E12").Select()
